# Applies the "Auto Doc Gen Spreadsheet" edit:
#  - B6 text changes from "$document: Document" to "document: Document"
#  - F6 text changes from "$document" to "document"
#  - the active selection moves from F6 to D11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "document: Document"
$ws.Range("F6").Value = "document"

$ws.Range("D11").Select()
